# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Cebolla" (Vega Monumental Concepción)
# right before the existing row 722, shifting the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 722-723; everything from old row 722 onward shifts down by 2.
$ws.Rows("722:723").Insert()

# Populate the newly inserted row 722
$ws.Range("A722").Value = 11
$ws.Range("B722").Value = "Vega Monumental Concepción"
$ws.Range("C722").Value = "Bíobío"
$ws.Range("D722").Value = 45077
$ws.Range("E722").Value = 8
$ws.Range("F722").Value = 100112004
$ws.Range("G722").Value = "Cebolla"
$ws.Range("H722").Value = "Sin especificar"
$ws.Range("I722").Value = "1a (guarda)"
$ws.Range("J722").Value = 600
$ws.Range("K722").Value = 8500
$ws.Range("L722").Value = 9000
$ws.Range("M722").Value = 8750
$ws.Range("N722").Value = "`$/malla 18 kilos"
$ws.Range("O722").Value = "Región Metropolitana"
$ws.Range("P722").Value = 486
$ws.Range("Q722").Value = 18
$ws.Range("R722").Value = "Hortaliza"

# Populate the newly inserted row 723
$ws.Range("A723").Value = 11
$ws.Range("B723").Value = "Vega Monumental Concepción"
$ws.Range("C723").Value = "Bíobío"
$ws.Range("D723").Value = 45077
$ws.Range("E723").Value = 8
$ws.Range("F723").Value = 100112004
$ws.Range("G723").Value = "Cebolla"
$ws.Range("H723").Value = "Sin especificar"
$ws.Range("I723").Value = "2a (guarda)"
$ws.Range("J723").Value = 300
$ws.Range("K723").Value = 7500
$ws.Range("L723").Value = 7500
$ws.Range("M723").Value = 7500
$ws.Range("N723").Value = "`$/malla 18 kilos"
$ws.Range("O723").Value = "Región Metropolitana"
$ws.Range("P723").Value = 417
$ws.Range("Q723").Value = 18
$ws.Range("R723").Value = "Hortaliza"
